$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 198. This shifts existing rows 198:291 down to 199:292,
# carrying their formatting with them.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical across the whole block (same market/category),
# so copy them straight from the row above (row 197) which is unaffected by the insert.
$ws.Cells.Item(198, 1).Value = $ws.Cells.Item(197, 1).Value2   # A: Mercado ID
$ws.Cells.Item(198, 2).Value = $ws.Cells.Item(197, 2).Value2   # B: Mercado
$ws.Cells.Item(198, 3).Value = $ws.Cells.Item(197, 3).Value2   # C: Region
$ws.Cells.Item(198, 4).Value = 44813                            # D: Fecha
$ws.Cells.Item(198, 4).NumberFormat = $ws.Cells.Item(197, 4).NumberFormat
$ws.Cells.Item(198, 5).Value = $ws.Cells.Item(197, 5).Value2   # E: Codreg
$ws.Cells.Item(198, 6).Value = $ws.Cells.Item(197, 6).Value2   # F: Categoria ID
$ws.Cells.Item(198, 7).Value = $ws.Cells.Item(197, 7).Value2   # G: Categoria
$ws.Cells.Item(198, 8).Value = $ws.Cells.Item(197, 8).Value2   # H: Variedad
$ws.Cells.Item(198, 9).Value = $ws.Cells.Item(197, 9).Value2   # I: Calidad
$ws.Cells.Item(198, 10).Value = 200                             # J: Volumen
$ws.Cells.Item(198, 11).Value = 28000                           # K: Precio minimo
$ws.Cells.Item(198, 12).Value = 30000                           # L: Precio maximo
$ws.Cells.Item(198, 13).Value = 29000                           # M: Precio promedio ponderado
$ws.Cells.Item(198, 14).Value = "$/caja 50 unidades"            # N: Unidad de comercializacion
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"  # O: Origen
$ws.Cells.Item(198, 16).Value = 580                             # P: Precio $/Kg
$ws.Cells.Item(198, 17).Value = 50                              # Q: Kg o Unidades
$ws.Cells.Item(198, 18).Value = $ws.Cells.Item(197, 18).Value2  # R: Clasificacion
